# Simulated Wild Card round and logged it
# Update the "R" (road/playoff) row totals on both the OFF and DEF sheets
# to reflect the simulated Wild Card round game.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$ws = $wb.Worksheets.Item("OFF")
$ws.Range("B3").Value = 193
$ws.Range("C3").Value = 126
$ws.Range("D3").Value = 39
$ws.Range("E3").Value = 23
$ws.Range("G3").Value = 1

# --- DEF sheet ---
$ws = $wb.Worksheets.Item("DEF")
$ws.Range("B3").Value = 218
$ws.Range("C3").Value = 148
$ws.Range("D3").Value = 42
$ws.Range("E3").Value = 16
